$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 38
$ws_ALC.Range("H38").Value = 307.75
$ws_ALC.Range("I38").Value = 72.4
$ws_ALC.Range("J38").Value = 475.85715
$ws_ALC.Range("K38").Value = 217.2
$ws_ALC.Range("L38").Value = 1427.57145
$ws_ALC.Range("M38").Value = 154.8
$ws_ALC.Range("N38").Value = -2171.57145

# ALC row 40
$ws_ALC.Range("H40").Value = 2232.0667
$ws_ALC.Range("J40").Value = 2112.625
$ws_ALC.Range("L40").Value = 2112.625
$ws_ALC.Range("N40").Value = -2462.625

# ALC row 42
$ws_ALC.Range("H42").Value = 129.6
$ws_ALC.Range("I42").Value = 66
$ws_ALC.Range("J42").Value = 225
$ws_ALC.Range("K42").Value = 198
$ws_ALC.Range("L42").Value = 675
$ws_ALC.Range("M42").Value = 32
$ws_ALC.Range("N42").Value = -1135

# ALC row 138
$ws_ALC.Range("H138").Value = 8734521
$ws_ALC.Range("I138").Value = 3463676.5
$ws_ALC.Range("J138").Value = 10640996
$ws_ALC.Range("K138").Value = 10391029.5
$ws_ALC.Range("L138").Value = 31922988
$ws_ALC.Range("M138").Value = -10385889.5
$ws_ALC.Range("N138").Value = -31933268

# ARM row 132
$ws_ARM.Range("H132").Value = 3256.7026
$ws_ARM.Range("I132").Value = 2787.2917
$ws_ARM.Range("J132").Value = 4123.3076
$ws_ARM.Range("K132").Value = 8361.875100000001
$ws_ARM.Range("L132").Value = 12369.9228
$ws_ARM.Range("M132").Value = -5831.875100000001
$ws_ARM.Range("N132").Value = -17429.9228

# BSM row 26
$ws_BSM.Range("H26").Value = 35975.4
$ws_BSM.Range("J26").Value = 49450
$ws_BSM.Range("L26").Value = 49450
$ws_BSM.Range("N26").Value = -50034

# BSM row 132
$ws_BSM.Range("H132").Value = 37223.637
$ws_BSM.Range("J132").Value = 37223.637
$ws_BSM.Range("L132").Value = 37223.637
$ws_BSM.Range("N132").Value = -47343.637

# BSM row 134
$ws_BSM.Range("H134").Value = 4243.59
$ws_BSM.Range("I134").Value = 2739.9167
$ws_BSM.Range("K134").Value = 8219.750100000001
$ws_BSM.Range("M134").Value = -5684.750100000001

# CRP row 22
$ws_CRP.Range("H22").Value = 522.2778
$ws_CRP.Range("J22").Value = 966.6667
$ws_CRP.Range("L22").Value = 966.6667
$ws_CRP.Range("N22").Value = -1666.6667

# CRP row 122
$ws_CRP.Range("H122").Value = 1382.7
$ws_CRP.Range("I122").Value = 1304.5
$ws_CRP.Range("K122").Value = 3913.5
$ws_CRP.Range("M122").Value = -1463.5

# CUL row 2
$ws_CUL.Range("H2").Value = 220.94444
$ws_CUL.Range("I2").Value = 143.33333
$ws_CUL.Range("J2").Value = 376.16666
$ws_CUL.Range("K2").Value = 859.9999799999999
$ws_CUL.Range("L2").Value = 2256.99996
$ws_CUL.Range("M2").Value = -746.9999799999999
$ws_CUL.Range("N2").Value = -2482.99996

# CUL row 38
$ws_CUL.Range("I38").Value = 157.14285
$ws_CUL.Range("J38").Value = 155.1
$ws_CUL.Range("K38").Value = 471.42855
$ws_CUL.Range("L38").Value = 465.3
$ws_CUL.Range("M38").Value = -124.42855
$ws_CUL.Range("N38").Value = -1159.3

# CUL row 41
$ws_CUL.Range("H41").Value = 1500
$ws_CUL.Range("I41").Value = 1500
$ws_CUL.Range("J41").Value = 0
$ws_CUL.Range("K41").Value = 4500
$ws_CUL.Range("L41").Value = 0
$ws_CUL.Range("N41").Value = ""
$ws_CUL.Range("M41").Value = -4162

# CUL row 87
$ws_CUL.Range("H87").Value = 16331.889
$ws_CUL.Range("I87").Value = 3495
$ws_CUL.Range("J87").Value = 19999.572
$ws_CUL.Range("K87").Value = 10485
$ws_CUL.Range("L87").Value = 59998.716
$ws_CUL.Range("M87").Value = -9237
$ws_CUL.Range("N87").Value = -62494.716

# CUL row 88
$ws_CUL.Range("H88").Value = 5000
$ws_CUL.Range("J88").Value = 5000
$ws_CUL.Range("L88").Value = 15000
$ws_CUL.Range("N88").Value = -15856

# CUL row 90
$ws_CUL.Range("H90").Value = 16331.889
$ws_CUL.Range("I90").Value = 3495
$ws_CUL.Range("J90").Value = 19999.572
$ws_CUL.Range("K90").Value = 31455
$ws_CUL.Range("L90").Value = 179996.148
$ws_CUL.Range("M90").Value = -25215
$ws_CUL.Range("N90").Value = -192476.148

# CUL row 91
$ws_CUL.Range("H91").Value = 5000
$ws_CUL.Range("J91").Value = 5000
$ws_CUL.Range("L91").Value = 15000
$ws_CUL.Range("N91").Value = -17964

# GSM row 80
$ws_GSM.Range("H80").Value = 3028.0557
$ws_GSM.Range("I80").Value = 2992.6924
$ws_GSM.Range("J80").Value = 3120
$ws_GSM.Range("K80").Value = 2992.6924
$ws_GSM.Range("L80").Value = 3120
$ws_GSM.Range("M80").Value = -1994.6924
$ws_GSM.Range("N80").Value = -5116

# GSM row 83
$ws_GSM.Range("H83").Value = 3028.0557
$ws_GSM.Range("I83").Value = 2992.6924
$ws_GSM.Range("J83").Value = 3120
$ws_GSM.Range("K83").Value = 14963.462
$ws_GSM.Range("L83").Value = 15600
$ws_GSM.Range("M83").Value = -9971.462
$ws_GSM.Range("N83").Value = -25584

# GSM row 107
$ws_GSM.Range("H107").Value = 1133.1875
$ws_GSM.Range("I107").Value = 1212.25
$ws_GSM.Range("J107").Value = 1054.125
$ws_GSM.Range("K107").Value = 1212.25
$ws_GSM.Range("L107").Value = 1054.125
$ws_GSM.Range("M107").Value = 707.75
$ws_GSM.Range("N107").Value = -4894.125

# GSM row 138
$ws_GSM.Range("H138").Value = 63900
$ws_GSM.Range("J138").Value = 63900
$ws_GSM.Range("L138").Value = 63900
$ws_GSM.Range("N138").Value = -74180

# GSM row 139
$ws_GSM.Range("H139").Value = 42500
$ws_GSM.Range("J139").Value = 42500
$ws_GSM.Range("L139").Value = 42500
$ws_GSM.Range("N139").Value = -52780

# LTW row 16
$ws_LTW.Range("H16").Value = 16670840
$ws_LTW.Range("I16").Value = 20003808
$ws_LTW.Range("J16").Value = 6000
$ws_LTW.Range("K16").Value = 20003808
$ws_LTW.Range("L16").Value = 6000
$ws_LTW.Range("M16").Value = -20003638
$ws_LTW.Range("N16").Value = -6340

# LTW row 86
$ws_LTW.Range("H86").Value = 40000
$ws_LTW.Range("J86").Value = 40000
$ws_LTW.Range("L86").Value = 40000
$ws_LTW.Range("N86").Value = -42372

# LTW row 89
$ws_LTW.Range("H89").Value = 40000
$ws_LTW.Range("J89").Value = 40000
$ws_LTW.Range("L89").Value = 120000
$ws_LTW.Range("N89").Value = -131856

# WVR row 41
$ws_WVR.Range("H41").Value = 6239.2
$ws_WVR.Range("I41").Value = 5342
$ws_WVR.Range("J41").Value = 6463.5
$ws_WVR.Range("K41").Value = 5342
$ws_WVR.Range("L41").Value = 6463.5
$ws_WVR.Range("N41").Value = -7243.5
$ws_WVR.Range("M41").Value = -4952

# WVR row 57
$ws_WVR.Range("H57").Value = 40796.668
$ws_WVR.Range("J57").Value = 33695
$ws_WVR.Range("L57").Value = 33695
$ws_WVR.Range("N57").Value = -35203

# WVR row 81
$ws_WVR.Range("H81").Value = 3646.8462
$ws_WVR.Range("I81").Value = 1405.75
$ws_WVR.Range("K81").Value = 2811.5
$ws_WVR.Range("M81").Value = -1750.5

# WVR row 84
$ws_WVR.Range("H84").Value = 3646.8462
$ws_WVR.Range("I84").Value = 1405.75
$ws_WVR.Range("K84").Value = 14057.5
$ws_WVR.Range("M84").Value = -8753.5

# WVR row 126
$ws_WVR.Range("H126").Value = 35246.656
$ws_WVR.Range("I126").Value = 50417.4
$ws_WVR.Range("J126").Value = 1533.8889
$ws_WVR.Range("K126").Value = 151252.2
$ws_WVR.Range("L126").Value = 4601.6667
$ws_WVR.Range("M126").Value = -148782.2
$ws_WVR.Range("N126").Value = -9308.3333

Write-Output "Applied Titan_Profits updates"